# Implementa a busca por similaridade.
# For a set of transactions that had no category (CATEGORIA, column F) yet,
# look at other transactions with similar/related items and copy over the
# category that was found for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$categorias = @{
    12  = "ALIMENTAÇÃO"
    21  = "ALIMENTAÇÃO"
    22  = "ALIMENTAÇÃO"
    75  = "ALIMENTAÇÃO"
    89  = "COMPRAS"
    92  = "ALIMENTAÇÃO"
    99  = "LAZER"
    119 = "TRANSPORTE"
    123 = "TRANSPORTE"
    126 = "MERCADO"
    129 = "CARRO"
    130 = "ALIMENTAÇÃO"
}

foreach ($row in $categorias.Keys) {
    $ws.Range("F$row").Value = $categorias[$row]
}
